$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy the formatting (bold, centered,
# bordered) from the existing header cell H1 so they match the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I (I0) and J (IF), rows 2-26
$data = @(
    @(2, 6, 7),
    @(3, 9, 9),
    @(4, 8, 8),
    @(5, 7, 7),
    @(6, 9, 9),
    @(7, 8, 8),
    @(8, 7, 7),
    @(9, 7, 8),
    @(10, 8, 8),
    @(11, 1, 2),
    @(12, 8, 9),
    @(13, 8, 8),
    @(14, 8, 8),
    @(15, 5, 6),
    @(16, 11, 11),
    @(17, 6, 7),
    @(18, 8, 8),
    @(19, 6, 7),
    @(20, 6, 7),
    @(21, 6, 7),
    @(22, 5, 6),
    @(23, 8, 8),
    @(24, 7, 7),
    @(25, 6, 6),
    @(26, 8, 8)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
